$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Insert a new row above the old row 13 ("Worksheet external file") to hold
# a new "Worksheet description" field, shifting everything below down by one.
$ws.Rows.Item(13).Insert() | Out-Null
$ws.Range("A13").Value = "Worksheet description"

# The two hyperlinked cells (originally B17 and B21) have shifted down to
# B18 and B22 along with the row insert above, but this environment does not
# automatically relocate existing Hyperlink objects with the row shift, so
# rebuild the hyperlinks collection to point at the correct, now-shifted
# cells (the mailto link on B8 is unaffected by the insert).
$ws.Range("A1").Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:d.orme@imperial.ac.uk") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B18"), "https://doi.org/10.1098/rstb.2011.0049") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B22"), "https://gtr.ukri.org/projects?ref=NE%2FK006339%2F1") | Out-Null

# Adding a hyperlink re-applies the built-in Hyperlink cell style but as a
# freshly duplicated style entry; reassign the canonical named style so the
# cells keep using the workbook's existing "Hyperlink" style.
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B18").Style = "Hyperlink"
$ws.Range("B22").Style = "Hyperlink"

# Make the Summary worksheet the active tab/selection (it previously was
# external_2), and select the newly added description cell, matching the
# workbook being reopened with focus on the new row.
$ws.Activate() | Out-Null
$ws.Range("B13").Select() | Out-Null
